$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This script rewrites the per-row financial metrics in the IFRS consolidated
# statement table (rows 2-9, columns D:AJ) to the corrected figures, and clears
# out the cells that no longer carry a value for that row/column combination.

# --- Row 2 ---
$ws.Range("D2").Value = 225009
$ws.Range("E2").Value = 26548
$ws.Range("F2").Value = 26548
$ws.Range("G2").Value = 28676
$ws.Range("H2").Value = 21996
$ws.Range("I2").Value = 20811
$ws.Range("J2").Value = 1185
$ws.Range("K2").Value = 3380218
$ws.Range("L2").Value = 3075069
$ws.Range("M2").Value = 305149
$ws.Range("N2").Value = 291841
$ws.Range("O2").Value = 13308
$ws.Range("P2").Value = 26451
$ws.Range("Q2").Value = -20806
$ws.Range("R2").Value = 9608
$ws.Range("S2").Value = 6869
$ws.Range("T2").Value = 1821
$ws.Range("V2").Value = 498447
$ws.Range("W2").Value = 11.8
$ws.Range("X2").Value = 9.779999999999999
$ws.Range("Y2").Value = 7.34
$ws.Range("Z2").Value = 0.68
$ws.Range("AA2").Value = 1007.73
$ws.Range("AB2").Value = 1053.66
$ws.Range("AC2").Value = 4288
$ws.Range("AD2").Value = 10.37
$ws.Range("AE2").Value = 60136
$ws.Range("AF2").Value = 0.74
$ws.Range("AG2").Value = 950
$ws.Range("AH2").Value = 2.14
$ws.Range("AI2").Value = 24.62
$ws.Range("AJ2").Value = 474199587
$ws.Range("U2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 230860
$ws.Range("E3").Value = 29731
$ws.Range("F3").Value = 29731
$ws.Range("G3").Value = 31406
$ws.Range("H3").Value = 24460
$ws.Range("I3").Value = 23672
$ws.Range("J3").Value = 788
$ws.Range("K3").Value = 3705480
$ws.Range("L3").Value = 3387349
$ws.Range("M3").Value = 318131
$ws.Range("N3").Value = 308397
$ws.Range("O3").Value = 9734
$ws.Range("P3").Value = 26451
$ws.Range("Q3").Value = 29703
$ws.Range("R3").Value = -52886
$ws.Range("S3").Value = 13128
$ws.Range("T3").Value = 1248
$ws.Range("V3").Value = 555956
$ws.Range("W3").Value = 12.88
$ws.Range("X3").Value = 10.6
$ws.Range("Y3").Value = 7.89
$ws.Range("Z3").Value = 0.6899999999999999
$ws.Range("AA3").Value = 1064.77
$ws.Range("AB3").Value = 1102.74
$ws.Range("AC3").Value = 4878
$ws.Range("AD3").Value = 8.109999999999999
$ws.Range("AE3").Value = 63548
$ws.Range("AF3").Value = 0.62
$ws.Range("AG3").Value = 1200
$ws.Range("AH3").Value = 3.03
$ws.Range("AI3").Value = 26.66
$ws.Range("AJ3").Value = 474199587
$ws.Range("U3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").Value = 222880
$ws.Range("E4").Value = 31086
$ws.Range("F4").Value = 31086
$ws.Range("G4").Value = 31705
$ws.Range("H4").Value = 28249
$ws.Range("I4").Value = 27748
$ws.Range("J4").Value = 501
$ws.Range("K4").Value = 3956803
$ws.Range("L4").Value = 3639353
$ws.Range("M4").Value = 317450
$ws.Range("N4").Value = 311097
$ws.Range("O4").Value = 6353
$ws.Range("P4").Value = 26451
$ws.Range("Q4").Value = 37972
$ws.Range("R4").Value = -75939
$ws.Range("S4").Value = 48446
$ws.Range("T4").Value = 2521
$ws.Range("V4").Value = 603962
$ws.Range("W4").Value = 13.95
$ws.Range("X4").Value = 12.68
$ws.Range("Y4").Value = 8.960000000000001
$ws.Range("Z4").Value = 0.74
$ws.Range("AA4").Value = 1146.43
$ws.Range("AB4").Value = 1100.16
$ws.Range("AC4").Value = 5810
$ws.Range("AD4").Value = 7.79
$ws.Range("AE4").Value = 65605
$ws.Range("AF4").Value = 0.6899999999999999
$ws.Range("AG4").Value = 1450
$ws.Range("AH4").Value = 3.2
$ws.Range("AI4").Value = 24.78
$ws.Range("AJ4").Value = 474199587
$ws.Range("U4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").Value = 239892
$ws.Range("E5").Value = 38300
$ws.Range("F5").Value = 38300
$ws.Range("G5").Value = 37976
$ws.Range("H5").Value = 29492
$ws.Range("I5").Value = 29188
$ws.Range("J5").Value = 304
$ws.Range("K5").Value = 4263070
$ws.Range("L5").Value = 3926034
$ws.Range("M5").Value = 337036
$ws.Range("N5").Value = 328202
$ws.Range("O5").Value = 8834
$ws.Range("P5").Value = 26451
$ws.Range("Q5").Value = 10213
$ws.Range("R5").Value = -107036
$ws.Range("S5").Value = 103325
$ws.Range("T5").Value = 1552
$ws.Range("V5").Value = 690009
$ws.Range("W5").Value = 15.97
$ws.Range("X5").Value = 12.29
$ws.Range("Y5").Value = 9.130000000000001
$ws.Range("Z5").Value = 0.72
$ws.Range("AA5").Value = 1164.87
$ws.Range("AB5").Value = 1174.21
$ws.Range("AC5").Value = 6155
$ws.Range("AD5").Value = 8.029999999999999
$ws.Range("AE5").Value = 69212
$ws.Range("AF5").Value = 0.71
$ws.Range("AG5").Value = 1450
$ws.Range("AH5").Value = 2.94
$ws.Range("AI5").Value = 23.56
$ws.Range("AJ5").Value = 474199587
$ws.Range("U5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = 238119
$ws.Range("E6").Value = 44994
$ws.Range("F6").Value = 44994
$ws.Range("G6").Value = 44666
$ws.Range("H6").Value = 31983
$ws.Range("I6").Value = 31567
$ws.Range("K6").Value = 4596005
$ws.Range("L6").Value = 4229491
$ws.Range("M6").Value = 366514
$ws.Range("N6").Value = 357256
$ws.Range("P6").Value = 26451
$ws.Range("Q6").Value = -63202
$ws.Range("R6").Value = -55126
$ws.Range("S6").Value = 138065
$ws.Range("T6").Value = 1429
$ws.Range("V6").Value = 839937
$ws.Range("W6").Value = 18.9
$ws.Range("X6").Value = 13.43
$ws.Range("Y6").Value = 9.210000000000001
$ws.Range("Z6").Value = 0.72
$ws.Range("AA6").Value = 1153.98
$ws.Range("AB6").Value = 1291.55
$ws.Range("AC6").Value = 6657
$ws.Range("AD6").Value = 5.95
$ws.Range("AE6").Value = 75923
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 1600
$ws.Range("AH6").Value = 4.04
$ws.Range("AI6").Value = 23.86
$ws.Range("AJ6").Value = 474199587
$ws.Range("U6").ClearContents()

# --- Row 7 ---
$ws.Range("E7").Value = 51422
$ws.Range("G7").Value = 50261
$ws.Range("H7").Value = 37062
$ws.Range("I7").Value = 34832
$ws.Range("K7").Value = 5379651
$ws.Range("L7").Value = 4960860
$ws.Range("M7").Value = 418789
$ws.Range("N7").Value = 391209
$ws.Range("P7").Value = 27157
$ws.Range("Y7").Value = 9.31
$ws.Range("Z7").Value = 0.74
$ws.Range("AA7").Value = 1184.57
$ws.Range("AC7").Value = 7167
$ws.Range("AD7").Value = 5.46
$ws.Range("AE7").Value = 81877
$ws.Range("AF7").Value = 0.48
$ws.Range("AG7").Value = 1803
$ws.Range("AH7").Value = 4.6
$ws.Range("AI7").Value = 24.54
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# --- Row 8 ---
$ws.Range("E8").Value = 51457
$ws.Range("G8").Value = 50370
$ws.Range("H8").Value = 36904
$ws.Range("I8").Value = 34531
$ws.Range("K8").Value = 5590714
$ws.Range("L8").Value = 5145505
$ws.Range("M8").Value = 445289
$ws.Range("N8").Value = 420728
$ws.Range("P8").Value = 27276
$ws.Range("Y8").Value = 8.51
$ws.Range("Z8").Value = 0.67
$ws.Range("AA8").Value = 1155.54
$ws.Range("AC8").Value = 7023
$ws.Range("AD8").Value = 5.57
$ws.Range("AE8").Value = 85569
$ws.Range("AF8").Value = 0.46
$ws.Range("AG8").Value = 1857
$ws.Range("AH8").Value = 4.74
$ws.Range("AI8").Value = 25.51
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# --- Row 9 ---
$ws.Range("E9").Value = 53549
$ws.Range("G9").Value = 52050
$ws.Range("H9").Value = 38023
$ws.Range("I9").Value = 35452
$ws.Range("K9").Value = 5768918
$ws.Range("L9").Value = 5291586
$ws.Range("M9").Value = 477332
$ws.Range("N9").Value = 452658
$ws.Range("P9").Value = 27264
$ws.Range("Y9").Value = 8.119999999999999
$ws.Range("Z9").Value = 0.67
$ws.Range("AA9").Value = 1108.58
$ws.Range("AC9").Value = 7210
$ws.Range("AD9").Value = 5.43
$ws.Range("AE9").Value = 92063
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 1964
$ws.Range("AH9").Value = 5.02
$ws.Range("AI9").Value = 26.27
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
